# Apply the "bento profile" query updates to the startup worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New StatQuery text (used by CasesTab, SamplesTab and FilesTab rows) ---
$statQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE samp.summarized_sample_type IN ["Normal Tissue"]  
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# --- Updated CasesTab "query" (adds the Cohort column) ---
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
WHERE samp.summarized_sample_type IN ["Normal Tissue"] 
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@

# --- Updated FilesTab "query" (drops the trailing Study Code column) ---
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
WHERE samp.summarized_sample_type IN ["Normal Tissue"] 
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis  
'@

# StatQuery column (C) now uses the same simplified counts query on every row
$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# FilesTab (row 4) query drops the trailing Study Code line
$ws.Range("B4").Value = $filesQuery

# CasesTab (row 2) query gets the extra Cohort column
$ws.Range("B2").Value = $casesQuery

# Row heights shrink now that the cell text is shorter (wrap-text autofit no longer hits the cap)
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210

# Reset the scroll position back to the top-left corner and move the selection to F3
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F3").Select() | Out-Null
